# Scheduled runner update: refresh market-board derived columns
# (currentAveragePrice*, LevePrice*, LeveProfit*) for the affected Leve
# rows across the profession sheets, per the latest price snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 71.111115
$ws.Cells.Item(2, 9).Value = 42.5
$ws.Cells.Item(2, 10).Value = 300
$ws.Cells.Item(2, 11).Value = 42.5
$ws.Cells.Item(2, 12).Value = 300
$ws.Cells.Item(2, 13).Value = 70.5
$ws.Cells.Item(2, 14).Value = -526
$ws.Cells.Item(8, 8).Value = 464.34042
$ws.Cells.Item(8, 9).Value = 1863
$ws.Cells.Item(8, 10).Value = 297.83334
$ws.Cells.Item(8, 11).Value = 5589
$ws.Cells.Item(8, 12).Value = 893.5000200000001
$ws.Cells.Item(8, 13).Value = -5450
$ws.Cells.Item(8, 14).Value = -1171.50002
$ws.Cells.Item(48, 8).Value = 10000
$ws.Cells.Item(48, 10).Value = 10000
$ws.Cells.Item(48, 12).Value = 30000
$ws.Cells.Item(48, 14).Value = -30584
$ws.Cells.Item(56, 8).Value = 10000
$ws.Cells.Item(56, 10).Value = 10000
$ws.Cells.Item(56, 12).Value = 30000
$ws.Cells.Item(56, 14).Value = -31068
$ws.Cells.Item(64, 8).Value = 6336
$ws.Cells.Item(64, 10).Value = 6405.375
$ws.Cells.Item(64, 12).Value = 6405.375
$ws.Cells.Item(64, 14).Value = -6901.375
$ws.Cells.Item(67, 8).Value = 6336
$ws.Cells.Item(67, 10).Value = 6405.375
$ws.Cells.Item(67, 12).Value = 6405.375
$ws.Cells.Item(67, 14).Value = -8121.375
$ws.Cells.Item(70, 8).Value = 72919256
$ws.Cells.Item(70, 9).Value = 62500700
$ws.Cells.Item(70, 11).Value = 187502100
$ws.Cells.Item(70, 13).Value = -187501830
$ws.Cells.Item(73, 8).Value = 72919256
$ws.Cells.Item(73, 9).Value = 62500700
$ws.Cells.Item(73, 11).Value = 187502100
$ws.Cells.Item(73, 13).Value = -187501164
$ws.Cells.Item(106, 8).Value = 3321.5
$ws.Cells.Item(106, 9).Value = 5072
$ws.Cells.Item(106, 11).Value = 5072
$ws.Cells.Item(106, 13).Value = -4441
$ws.Cells.Item(132, 8).Value = 1983.5834
$ws.Cells.Item(132, 9).Value = 1590.45
$ws.Cells.Item(132, 10).Value = 3949.25
$ws.Cells.Item(132, 11).Value = 4771.35
$ws.Cells.Item(132, 12).Value = 11847.75
$ws.Cells.Item(132, 13).Value = -2241.35
$ws.Cells.Item(132, 14).Value = -16907.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(63, 8).Value = 2709.8
$ws.Cells.Item(63, 9).Value = 1137.25
$ws.Cells.Item(63, 11).Value = 1137.25
$ws.Cells.Item(63, 13).Value = -451.25
$ws.Cells.Item(66, 8).Value = 2709.8
$ws.Cells.Item(66, 9).Value = 1137.25
$ws.Cells.Item(66, 11).Value = 5686.25
$ws.Cells.Item(66, 13).Value = -2254.25
$ws.Cells.Item(74, 8).Value = 35003.71
$ws.Cells.Item(74, 9).Value = 49343.855
$ws.Cells.Item(74, 10).Value = 4889.4
$ws.Cells.Item(74, 11).Value = 49343.855
$ws.Cells.Item(74, 12).Value = 4889.4
$ws.Cells.Item(74, 13).Value = -48469.855
$ws.Cells.Item(74, 14).Value = -6637.4
$ws.Cells.Item(77, 8).Value = 35003.71
$ws.Cells.Item(77, 9).Value = 49343.855
$ws.Cells.Item(77, 10).Value = 4889.4
$ws.Cells.Item(77, 11).Value = 246719.275
$ws.Cells.Item(77, 12).Value = 24447
$ws.Cells.Item(77, 13).Value = -242351.275
$ws.Cells.Item(77, 14).Value = -33183
$ws.Cells.Item(101, 8).Value = 70000
$ws.Cells.Item(101, 10).Value = 70000
$ws.Cells.Item(101, 12).Value = 70000
$ws.Cells.Item(101, 14).Value = -76490
$ws.Cells.Item(102, 8).Value = 4757.1113
$ws.Cells.Item(102, 9).Value = 3830.5715
$ws.Cells.Item(102, 10).Value = 8000
$ws.Cells.Item(102, 11).Value = 3830.5715
$ws.Cells.Item(102, 12).Value = 8000
$ws.Cells.Item(102, 13).Value = -2208.5715
$ws.Cells.Item(102, 14).Value = -11244
$ws.Cells.Item(132, 8).Value = 5719.6763
$ws.Cells.Item(132, 9).Value = 3981.2
$ws.Cells.Item(132, 11).Value = 11943.6
$ws.Cells.Item(132, 13).Value = -9413.599999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 127436.375
$ws.Cells.Item(86, 9).Value = 168583.5
$ws.Cells.Item(86, 11).Value = 168583.5
$ws.Cells.Item(86, 13).Value = -167460.5
$ws.Cells.Item(89, 8).Value = 127436.375
$ws.Cells.Item(89, 9).Value = 168583.5
$ws.Cells.Item(89, 11).Value = 842917.5
$ws.Cells.Item(89, 13).Value = -837301.5
$ws.Cells.Item(94, 8).Value = 2661.1785
$ws.Cells.Item(94, 9).Value = 1277.909
$ws.Cells.Item(94, 11).Value = 1277.909
$ws.Cells.Item(94, 13).Value = -826.9090000000001
$ws.Cells.Item(105, 8).Value = 3386.0278
$ws.Cells.Item(105, 9).Value = 2609.4546
$ws.Cells.Item(105, 10).Value = 4606.357
$ws.Cells.Item(105, 11).Value = 2609.4546
$ws.Cells.Item(105, 12).Value = 4606.357
$ws.Cells.Item(105, 13).Value = -862.4546
$ws.Cells.Item(105, 14).Value = -8100.357
$ws.Cells.Item(118, 8).Value = 69444
$ws.Cells.Item(118, 10).Value = 69444
$ws.Cells.Item(118, 12).Value = 69444
$ws.Cells.Item(118, 14).Value = -72758
$ws.Cells.Item(134, 8).Value = 3734328.8
$ws.Cells.Item(134, 9).Value = 4466309.5
$ws.Cells.Item(134, 11).Value = 13398928.5
$ws.Cells.Item(134, 13).Value = -13396393.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 4593.375
$ws.Cells.Item(31, 9).Value = 1821.8334
$ws.Cells.Item(31, 10).Value = 6672.0312
$ws.Cells.Item(31, 11).Value = 1821.8334
$ws.Cells.Item(31, 12).Value = 6672.0312
$ws.Cells.Item(31, 13).Value = -1526.8334
$ws.Cells.Item(31, 14).Value = -7262.0312
$ws.Cells.Item(34, 8).Value = 4593.375
$ws.Cells.Item(34, 9).Value = 1821.8334
$ws.Cells.Item(34, 10).Value = 6672.0312
$ws.Cells.Item(34, 11).Value = 1821.8334
$ws.Cells.Item(34, 12).Value = 6672.0312
$ws.Cells.Item(34, 13).Value = -1619.8334
$ws.Cells.Item(34, 14).Value = -7076.0312
$ws.Cells.Item(94, 8).Value = 2099
$ws.Cells.Item(94, 10).Value = 1722.5714
$ws.Cells.Item(94, 12).Value = 1722.5714
$ws.Cells.Item(94, 14).Value = -2624.5714
$ws.Cells.Item(99, 8).Value = 5822.923
$ws.Cells.Item(99, 10).Value = 7924.875
$ws.Cells.Item(99, 12).Value = 7924.875
$ws.Cells.Item(99, 14).Value = -10920.875
$ws.Cells.Item(126, 8).Value = 5822.923
$ws.Cells.Item(126, 10).Value = 7924.875
$ws.Cells.Item(126, 12).Value = 23774.625
$ws.Cells.Item(126, 14).Value = -28714.625
$ws.Cells.Item(132, 8).Value = 4508.129
$ws.Cells.Item(132, 9).Value = 2777.3157
$ws.Cells.Item(132, 10).Value = 7248.5835
$ws.Cells.Item(132, 11).Value = 8331.947100000001
$ws.Cells.Item(132, 12).Value = 21745.7505
$ws.Cells.Item(132, 13).Value = -5801.947100000001
$ws.Cells.Item(132, 14).Value = -26805.7505
$ws.Cells.Item(134, 8).Value = 3386.7593
$ws.Cells.Item(134, 9).Value = 1483.6154
$ws.Cells.Item(134, 11).Value = 4450.8462
$ws.Cells.Item(134, 13).Value = -1915.8462

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(103, 8).Value = 956.75
$ws.Cells.Item(103, 9).Value = 110
$ws.Cells.Item(103, 10).Value = 1077.7142
$ws.Cells.Item(103, 11).Value = 330
$ws.Cells.Item(103, 12).Value = 3233.1426
$ws.Cells.Item(103, 13).Value = 549
$ws.Cells.Item(103, 14).Value = -4991.142599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 1283.8334
$ws.Cells.Item(46, 10).Value = 2018.75
$ws.Cells.Item(46, 12).Value = 2018.75
$ws.Cells.Item(46, 14).Value = -2394.75
$ws.Cells.Item(74, 8).Value = 14000
$ws.Cells.Item(74, 9).Value = 14000
$ws.Cells.Item(74, 10).Value = 0
$ws.Cells.Item(74, 11).Value = 14000
$ws.Cells.Item(74, 12).Value = 0
$ws.Cells.Item(74, 13).Value = -13002
$ws.Cells.Item(74, 14).Value = $null
$ws.Cells.Item(77, 8).Value = 14000
$ws.Cells.Item(77, 9).Value = 14000
$ws.Cells.Item(77, 10).Value = 0
$ws.Cells.Item(77, 11).Value = 42000
$ws.Cells.Item(77, 12).Value = 0
$ws.Cells.Item(77, 13).Value = -37008
$ws.Cells.Item(77, 14).Value = $null
$ws.Cells.Item(100, 8).Value = 4219.067
$ws.Cells.Item(100, 9).Value = 3598.875
$ws.Cells.Item(100, 10).Value = 4927.857
$ws.Cells.Item(100, 11).Value = 3598.875
$ws.Cells.Item(100, 12).Value = 4927.857
$ws.Cells.Item(100, 13).Value = -3057.875
$ws.Cells.Item(100, 14).Value = -6009.857
$ws.Cells.Item(112, 8).Value = 52387
$ws.Cells.Item(112, 10).Value = 52387
$ws.Cells.Item(112, 12).Value = 52387
$ws.Cells.Item(112, 14).Value = -55341
$ws.Cells.Item(122, 8).Value = 4327.778
$ws.Cells.Item(122, 9).Value = 2983.3333
$ws.Cells.Item(122, 10).Value = 5000
$ws.Cells.Item(122, 11).Value = 8949.999899999999
$ws.Cells.Item(122, 12).Value = 15000
$ws.Cells.Item(122, 13).Value = -6499.999899999999
$ws.Cells.Item(122, 14).Value = -19900
$ws.Cells.Item(136, 8).Value = 8124.7925
$ws.Cells.Item(136, 9).Value = 2533.818
$ws.Cells.Item(136, 10).Value = 17349.9
$ws.Cells.Item(136, 11).Value = 7601.454000000001
$ws.Cells.Item(136, 12).Value = 52049.7
$ws.Cells.Item(136, 13).Value = -5051.454000000001
$ws.Cells.Item(136, 14).Value = -57149.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 14).Value = $null
$ws.Cells.Item(132, 8).Value = 8061.0586
$ws.Cells.Item(132, 9).Value = 8233.23
$ws.Cells.Item(132, 11).Value = 24699.69
$ws.Cells.Item(132, 13).Value = -22169.69

